# cryptos.xlsx — scheduled scrape refresh (GitHub Actions, 2023-09-03 19:45 UTC).
# Updates Price (D) and Volume(1h) (E) for each coin row, and re-sorts three
# adjacent coin pairs whose ranking order flipped: Polkadot/WrappedEther
# (rows 12-13), TrustWalletToken/BabyDogeCoin (rows 42-43), and
# Frax/Aave (rows 48-49) — for those, Coin/Link/Price/Volume all move together.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new Price text parses as a plain number (e.g. "1.002") need to be
# forced to Text format first, otherwise Excel auto-converts them to a numeric
# value instead of keeping the literal string used by the source data feed.
$numericLookingCells = @(
    "D4",
    "D5",
    "D6",
    "D7",
    "D8",
    "D9",
    "D10",
    "D11",
    "D13",
    "D14",
    "D16",
    "D19",
    "D20",
    "D21",
    "D22",
    "D24",
    "D25",
    "D26",
    "D27",
    "D28",
    "D29",
    "D31",
    "D32",
    "D33",
    "D34",
    "D35",
    "D37",
    "D38",
    "D39",
    "D41",
    "D43",
    "D44",
    "D45",
    "D47",
    "D48",
    "D49",
    "D50",
    "D51"
)
foreach ($cellRef in $numericLookingCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = "26.019.18"
$ws.Range("E2").Value = "  +0.60%  "
# Row 3
$ws.Range("D3").Value = "1.640.87"
$ws.Range("E3").Value = "  +0.41%  "
# Row 4
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.09%  "
# Row 5
$ws.Range("D5").Value = "214.96"
$ws.Range("E5").Value = "  +0.52%  "
# Row 6
$ws.Range("D6").Value = "0.5090"
$ws.Range("E6").Value = "  +1.46%  "
# Row 7
$ws.Range("D7").Value = "1.002"
$ws.Range("E7").Value = "  +0.25%  "
# Row 8
$ws.Range("D8").Value = "0.2567"
$ws.Range("E8").Value = "  +0.27%  "
# Row 9
$ws.Range("D9").Value = "0.06376"
$ws.Range("E9").Value = "  +0.30%  "
# Row 10
$ws.Range("D10").Value = "19.52"
$ws.Range("E10").Value = "  +0.57%  "
# Row 11
$ws.Range("D11").Value = "0.07768"
$ws.Range("E11").Value = "  +0.06%  "
# Row 12
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.673.39"
$ws.Range("E12").Value = "  +2.42%  "
# Row 13
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "4.288"
$ws.Range("E13").Value = "  +1.02%  "
# Row 14
$ws.Range("D14").Value = "0.5449"
$ws.Range("E14").Value = "  +0.99%  "
# Row 15
$ws.Range("D15").Value = "0.0₅7743"
$ws.Range("E15").Value = "  -1.30%  "
# Row 16
$ws.Range("D16").Value = "64.25"
$ws.Range("E16").Value = "  -0.10%  "
# Row 17
$ws.Range("D17").Value = "26.035.80"
$ws.Range("E17").Value = "  +0.70%  "
# Row 18
$ws.Range("E18").Value = "  +0.30%  "
# Row 19
$ws.Range("D19").Value = "196.65"
$ws.Range("E19").Value = "  -0.18%  "
# Row 20
$ws.Range("D20").Value = "4.426"
$ws.Range("E20").Value = "  +1.53%  "
# Row 21
$ws.Range("D21").Value = "9.931"
$ws.Range("E21").Value = "  +0.48%  "
# Row 22
$ws.Range("D22").Value = "6.044"
$ws.Range("E22").Value = "  +1.70%  "
# Row 23
$ws.Range("E23").Value = "  +0.31%  "
# Row 24
$ws.Range("D24").Value = "1.884"
$ws.Range("E24").Value = "  +0.21%  "
# Row 25
$ws.Range("D25").Value = "141.19"
$ws.Range("E25").Value = "  +1.22%  "
# Row 26
$ws.Range("D26").Value = "0.1196"
$ws.Range("E26").Value = "  +5.35%  "
# Row 27
$ws.Range("D27").Value = "6.840"
$ws.Range("E27").Value = "  +0.45%  "
# Row 28
$ws.Range("D28").Value = "15.58"
$ws.Range("E28").Value = "  -0.29%  "
# Row 29
$ws.Range("D29").Value = "1.236"
$ws.Range("E29").Value = "  +0.29%  "
# Row 30
$ws.Range("E30").Value = "  +0.44%  "
# Row 31
$ws.Range("D31").Value = "3.258"
$ws.Range("E31").Value = "  +0.42%  "
# Row 32
$ws.Range("D32").Value = "3.171"
$ws.Range("E32").Value = "  +0.08%  "
# Row 33
$ws.Range("D33").Value = "1.528"
$ws.Range("E33").Value = "  +0.25%  "
# Row 34
$ws.Range("D34").Value = "2.363"
$ws.Range("E34").Value = "  +0.49%  "
# Row 35
$ws.Range("D35").Value = "0.8942"
$ws.Range("E35").Value = "  +1.14%  "
# Row 36
$ws.Range("D36").Value = "1.144.30"
$ws.Range("E36").Value = "  +2.01%  "
# Row 37
$ws.Range("D37").Value = "2.582"
$ws.Range("E37").Value = "  -0.29%  "
# Row 38
$ws.Range("D38").Value = "0.5454"
$ws.Range("E38").Value = "  -0.84%  "
# Row 39
$ws.Range("D39").Value = "0.01557"
$ws.Range("E39").Value = "  +0.37%  "
# Row 40
$ws.Range("E40").Value = "  +0.34%  "
# Row 41
$ws.Range("D41").Value = "2.521"
$ws.Range("E41").Value = "  -1.07%  "
# Row 42
$ws.Range("B42").Value = "BabyDogeCoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D42").Value = "0.0₈128"
$ws.Range("E42").Value = "  +5.04%  "
# Row 43
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "0.8100"
$ws.Range("E43").Value = "  +0.14%  "
# Row 44
$ws.Range("D44").Value = "99.16"
$ws.Range("E44").Value = "  +0.07%  "
# Row 45
$ws.Range("D45").Value = "5.428"
$ws.Range("E45").Value = "  -4.05%  "
# Row 46
$ws.Range("D46").Value = "1.777.72"
$ws.Range("E46").Value = "  +0.37%  "
# Row 47
$ws.Range("D47").Value = "0.4527"
$ws.Range("E47").Value = "  +0.66%  "
# Row 48
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").Value = "54.98"
$ws.Range("E48").Value = "  -0.20%  "
# Row 49
$ws.Range("B49").Value = "Frax"
$ws.Range("C49").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D49").Value = "1.000"
$ws.Range("E49").Value = "  -0.21%  "
# Row 50
$ws.Range("D50").Value = "0.05056"
$ws.Range("E50").Value = "  -0.10%  "
# Row 51
$ws.Range("D51").Value = "1.003"
$ws.Range("E51").Value = "  +0.19%  "

# Restore default (General) formatting now that the text values are committed,
# so the cells end up style-free just like the rest of the data rows.
foreach ($cellRef in $numericLookingCells) {
    $ws.Range($cellRef).ClearFormats()
}
